# Youtube_trends_slides.pptx edit:
#   Slide 4 - the rhetorical question under "Childish Gambino is trending... A lot."
#   is reworded from
#       "But what categories are the most trending?"
#   to
#       "But what categories are trending the most?"
#
# The sentence lives in the same paragraph/run-less text box as the preceding
# sentence (separated only by a soft line break), so we must only touch the
# characters belonging to this second sentence and leave the first sentence's
# run (and its ellipsis character) completely untouched.

$p = $ppt.ActivePresentation

$oldSentence = "But what categories are the most trending?"
$newSentence = "But what categories are trending the most?"

$found = $false

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)

        if (-not $shape.HasTextFrame) {
            continue
        }

        $textRange = $shape.TextFrame.TextRange
        $fullText = $textRange.Text

        $charIndex = $fullText.IndexOf($oldSentence)
        if ($charIndex -ge 0) {
            # COM text ranges are 1-indexed; IndexOf() is 0-indexed.
            $sentenceRange = $textRange.Characters($charIndex + 1, $oldSentence.Length)
            $sentenceRange.Text = $newSentence
            $found = $true
        }
    }
}

if (-not $found) {
    throw "Could not locate the target sentence to update."
}
